$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-6, columns B (TB), C (d2S), D (K), E (IP), G (sum)
# Column A (date) and F (Win) are unchanged.

$data = @{
    2 = @{ B = 0.1554434735375247;  C = 0.0001537489499301437; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.368216780576502 }
    3 = @{ B = 0.06328177979961902; C = 0.0001537489499301437; D = 0.1529057820181812; E = 0.4998867070740569; G = 0.7162280178417872 }
    4 = @{ B = 3.182878228561681;   C = 1.65323645889881;      D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    5 = @{ B = 3.182878228561681;   C = 1.65323645889881;      D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    6 = @{ B = 3.182878228561681;   C = 1.65323645889881;      D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
}
